$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format column D as Text so numeric-looking price strings (e.g. "536.32")
# are preserved as literal text instead of being auto-converted to numbers,
# matching the inlineStr text cells in the target workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.088.94"
$ws.Range("E2").Value = "  +1.49%  "
$ws.Range("D3").Value = "3.136.90"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "536.32"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").Value = "139.04"
$ws.Range("E6").Value = "  +2.58%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.514"
$ws.Range("E8").Value = "  +10.67%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("D11").Value = "0.423"
$ws.Range("E11").Value = "  +5.50%  "
$ws.Range("E12").Value = "  +2.88%  "
$ws.Range("D13").Value = "3.679.71"
$ws.Range("E13").Value = "  +2.15%  "
$ws.Range("D14").Value = "25.95"
$ws.Range("E14").Value = "  +2.77%  "
$ws.Range("E15").Value = "  +4.98%  "
$ws.Range("D16").Value = "58.213.29"
$ws.Range("E16").Value = "  +1.61%  "
$ws.Range("D17").Value = "3.144.04"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("D18").Value = "6.23"
$ws.Range("E18").Value = "  +6.32%  "
$ws.Range("D19").Value = "12.97"
$ws.Range("E19").Value = "  +4.28%  "
$ws.Range("D20").Value = "8.17"
$ws.Range("E20").Value = "  +4.57%  "
$ws.Range("D21").Value = "375.03"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("E23").Value = "  -0.65%  "
$ws.Range("D24").Value = "70.36"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").Value = "0.515"
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0883"
$ws.Range("E28").Value = "  +2.22%  "
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").Value = "7.95"
$ws.Range("E29").Value = "  +10.42%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "1.89"
$ws.Range("E30").Value = "  +1.83%  "
$ws.Range("B31").Value = "RenderToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D31").Value = "6.15"
$ws.Range("E31").Value = "  +5.76%  "
$ws.Range("D32").Value = "21.75"
$ws.Range("E32").Value = "  +4.08%  "
$ws.Range("D33").Value = "5.15"
$ws.Range("E33").Value = "  +6.44%  "
$ws.Range("E34").Value = "  +3.77%  "
$ws.Range("D35").Value = "161.74"
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").Value = "6.26"
$ws.Range("E36").Value = "  +4.40%  "
$ws.Range("D37").Value = "1.36"
$ws.Range("E37").Value = "  +11.05%  "
$ws.Range("D38").Value = "25.53"
$ws.Range("E38").Value = "  +0.61%  "
$ws.Range("E39").Value = "  +6.54%  "
$ws.Range("D40").Value = "2.643.83"
$ws.Range("E40").Value = "  +9.96%  "
$ws.Range("D41").Value = "0.0681"
$ws.Range("E41").Value = "  +3.83%  "
$ws.Range("D42").Value = "4.25"
$ws.Range("E42").Value = "  +4.86%  "
$ws.Range("D43").Value = "38.43"
$ws.Range("E43").Value = "  +5.00%  "
$ws.Range("D44").Value = "0.700"
$ws.Range("E44").Value = "  +1.20%  "
$ws.Range("E45").Value = "  +4.71%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  +12.04%  "
$ws.Range("D48").Value = "6.23"
$ws.Range("E48").Value = "  +4.30%  "
$ws.Range("D49").Value = "0.977"
$ws.Range("E49").Value = "  +3.60%  "
$ws.Range("D50").Value = "20.26"
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("D51").Value = "0.749"
$ws.Range("E51").Value = "  -0.23%  "
